# Update the "Date" column (B) of the NumberError RAD test-result sheet
# with the timestamps from the latest test run (Mon Oct 02 2023), as
# captured by the new Katalon RAD test execution for MRF.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDates = @(
    "Mon Oct 02 16:33:43 EDT 2023",
    "Mon Oct 02 16:33:52 EDT 2023",
    "Mon Oct 02 16:34:02 EDT 2023",
    "Mon Oct 02 16:34:12 EDT 2023",
    "Mon Oct 02 16:34:22 EDT 2023",
    "Mon Oct 02 16:34:32 EDT 2023",
    "Mon Oct 02 16:34:41 EDT 2023",
    "Mon Oct 02 16:34:51 EDT 2023",
    "Mon Oct 02 16:35:01 EDT 2023",
    "Mon Oct 02 16:35:11 EDT 2023",
    "Mon Oct 02 16:35:20 EDT 2023",
    "Mon Oct 02 16:35:30 EDT 2023",
    "Mon Oct 02 16:35:40 EDT 2023",
    "Mon Oct 02 16:35:49 EDT 2023",
    "Mon Oct 02 16:35:59 EDT 2023",
    "Mon Oct 02 16:36:09 EDT 2023"
)

# Rows 2..17 hold the 16 test-result records; column B is "Date".
for ($i = 0; $i -lt $newDates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $newDates[$i]
}
